# The weekly refresh re-sorted the daily Chirimoya price rows (D=Fecha,
# L=Calidad, M=Volumen, N=Precio minimo, O=Precio maximo,
# P=Precio promedio ponderado, S=Precio $/Kg) for rows 3,4,7,8,9,10.
# Apply the new values cell by cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("D3").Value = 44452
$ws.Range("L3").Value = "Primera"
$ws.Range("N3").Value = 21000
$ws.Range("O3").Value = 22000
$ws.Range("P3").Value = 21500
$ws.Range("S3").Value = 2150

# Row 4
$ws.Range("D4").Value = 44447
$ws.Range("M4").Value = 60
$ws.Range("N4").Value = 21000
$ws.Range("O4").Value = 22000
$ws.Range("P4").Value = 21500
$ws.Range("S4").Value = 2150

# Row 7
$ws.Range("D7").Value = 44487
$ws.Range("M7").Value = 30
$ws.Range("N7").Value = 23000
$ws.Range("O7").Value = 24000
$ws.Range("P7").Value = 23500
$ws.Range("S7").Value = 2350

# Row 8
$ws.Range("D8").Value = 44460
$ws.Range("L8").Value = "Especial"
$ws.Range("N8").Value = 31000
$ws.Range("O8").Value = 32000
$ws.Range("P8").Value = 31500
$ws.Range("S8").Value = 3150

# Row 9
$ws.Range("D9").Value = 44460
$ws.Range("N9").Value = 30000
$ws.Range("O9").Value = 30000
$ws.Range("P9").Value = 30000
$ws.Range("S9").Value = 3000

# Row 10
$ws.Range("D10").Value = 44448
